$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UILocalizations")

# Insert a new row at row 30, shifting existing rows down.
$ws.Rows.Item(30).Insert()
$ws.Rows.Item(30).RowHeight = 15.75

# Fill in the new row's values. Shared-string insertion order: key, then
# Dutch (nl), then English (en), to mirror the author's save order.
$ws.Cells.Item(30, 1).Value = "VOICE_PRESS_TO_TALK"
$ws.Cells.Item(30, 3).Value = "Druk hier om te praten"
$ws.Cells.Item(30, 2).Value = "Press here to talk"

# Update the view to match (topLeftCell A7, selection A30).
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("A30").Select()
